$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.891.13"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.368.13"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'569.78"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'138.58"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'7.67"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("E11").Value = "  -5.02%  "
$ws.Range("D12").Value = "3.944.62"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "'27.69"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "3.374.48"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "60.974.39"
$ws.Range("E18").Value = "  -3.05%  "
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("D20").Value = "'8.88"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "'380.62"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'75.51"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -5.68%  "
$ws.Range("E26").Value = "  +6.72%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'7.15"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'1.33"
$ws.Range("E32").Value = "  -7.09%  "
$ws.Range("D33").Value = "'22.89"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "'6.85"
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("D35").Value = "'167.28"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").Value = "'4.91"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "3.403.67"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").Value = "'0.0759"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").Value = "'25.34"
$ws.Range("E40").Value = "  -9.25%  "
$ws.Range("D41").Value = "'0.772"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "'4.32"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").Value = "2.454.11"
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("D48").Value = "'22.19"
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("D50").Value = "'2.01"
$ws.Range("E51").Value = "  -3.75%  "
